$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace underscores with spaces in the "Employee Grade" column (G) values like LEVEL_5 -> LEVEL 5
for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $text = $cell.Text
    if ($text -ne $null -and $text -match '^LEVEL_\d+$') {
        $cell.Value = $text -replace '_', ' '
    }
}
